# Fri, Jun 12, 2020  4:05:17 PM
#
# 1) Re-point the table on slide 5 at the "Office Theme"-family built-in
#    table style (Medium Style 2 - Accent 1) instead of the custom
#    "Table_0" style that ships in ppt/tableStyles.xml.
# 2) Swap the presentation's colour scheme over to the stock "Office"
#    palette (the deck's theme1.xml/theme2.xml pair effectively traded
#    places: the design actually applied to the slide master becomes the
#    "Office Theme" colours).

$p = $ppt.ActivePresentation

# --- 1) table style -------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$tableShape = $slide5.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{F75068C4-E26D-4801-B351-4DC757F414C3}", $false)

# --- 2) theme colours -------------------------------------------------------
function Hex2BGR([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $tcs.Item($i).RGB = Hex2BGR($officeColors[$i - 1])
}
